$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2 (hunk 0)
$ws.Cells.Item(2, 8).Value = 75
$ws.Cells.Item(2, 9).Value = 75
$ws.Cells.Item(2, 11).Value = 75
$ws.Cells.Item(2, 13).Value = 38

# row 18 (hunk 1)
$ws.Cells.Item(18, 8).Value = 876.5
$ws.Cells.Item(18, 9).Value = 876.5
$ws.Cells.Item(18, 11).Value = 876.5
$ws.Cells.Item(18, 13).Value = -592.5

# row 76 (hunk 2)
$ws.Cells.Item(76, 8).Value = 2987.5
$ws.Cells.Item(76, 9).Value = 2987.5
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 2987.5
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -2672.5
$ws.Cells.Item(76, 14).ClearContents()

# row 79 (hunk 3)
$ws.Cells.Item(79, 8).Value = 2987.5
$ws.Cells.Item(79, 9).Value = 2987.5
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 2987.5
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -1895.5
$ws.Cells.Item(79, 14).ClearContents()

# row 107 (hunk 4)
$ws.Cells.Item(107, 8).Value = 444
$ws.Cells.Item(107, 10).Value = 1005
$ws.Cells.Item(107, 12).Value = 1005
$ws.Cells.Item(107, 14).Value = -4845

# row 129 (hunk 5)
$ws.Cells.Item(129, 8).Value = 1233.8572
$ws.Cells.Item(129, 9).Value = 1022.8333
$ws.Cells.Item(129, 11).Value = 3068.4999
$ws.Cells.Item(129, 13).Value = 1931.5001

# row 132 (hunk 6)
$ws.Cells.Item(132, 8).Value = 9137.25
$ws.Cells.Item(132, 9).Value = 14524.5
$ws.Cells.Item(132, 11).Value = 43573.5
$ws.Cells.Item(132, 13).Value = -41043.5

# row 138 (hunk 7)
$ws.Cells.Item(138, 8).Value = 1980
$ws.Cells.Item(138, 10).Value = 2800
$ws.Cells.Item(138, 12).Value = 8400
$ws.Cells.Item(138, 14).Value = -18680


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 28 (hunk 8)
$ws.Cells.Item(28, 8).Value = 35000
$ws.Cells.Item(28, 9).Value = 35000
$ws.Cells.Item(28, 11).Value = 35000
$ws.Cells.Item(28, 13).Value = -34808

# row 32 (hunk 9)
$ws.Cells.Item(32, 8).Value = 17681.285
$ws.Cells.Item(32, 9).Value = 14753.8
$ws.Cells.Item(32, 11).Value = 14753.8
$ws.Cells.Item(32, 13).Value = -14466.8

# row 88 (hunk 10)
$ws.Cells.Item(88, 8).Value = 1639.6
$ws.Cells.Item(88, 10).Value = 1583.3334
$ws.Cells.Item(88, 12).Value = 1583.3334
$ws.Cells.Item(88, 14).Value = -2395.3334

# row 91 (hunk 11)
$ws.Cells.Item(91, 8).Value = 1639.6
$ws.Cells.Item(91, 10).Value = 1583.3334
$ws.Cells.Item(91, 12).Value = 1583.3334
$ws.Cells.Item(91, 14).Value = -4391.3334

# row 92 (hunk 12)
$ws.Cells.Item(92, 8).Value = 62999.4
$ws.Cells.Item(92, 10).Value = 56249.25
$ws.Cells.Item(92, 12).Value = 56249.25
$ws.Cells.Item(92, 14).Value = -61241.25

# row 99 (hunk 13)
$ws.Cells.Item(99, 8).Value = 35000
$ws.Cells.Item(99, 9).Value = 35000
$ws.Cells.Item(99, 11).Value = 35000
$ws.Cells.Item(99, 13).Value = -32005

# row 132 (hunk 14)
$ws.Cells.Item(132, 8).Value = 3267.3333
$ws.Cells.Item(132, 9).Value = 2425.75
$ws.Cells.Item(132, 11).Value = 7277.25
$ws.Cells.Item(132, 13).Value = -4747.25


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 29 (hunk 15)
$ws.Cells.Item(29, 8).Value = 29000
$ws.Cells.Item(29, 9).Value = 29000
$ws.Cells.Item(29, 11).Value = 29000
$ws.Cells.Item(29, 13).Value = -28711

# row 86 (hunk 16)
$ws.Cells.Item(86, 8).Value = 1392.3
$ws.Cells.Item(86, 9).Value = 1502.1428
$ws.Cells.Item(86, 10).Value = 1136
$ws.Cells.Item(86, 11).Value = 1502.1428
$ws.Cells.Item(86, 12).Value = 1136
$ws.Cells.Item(86, 13).Value = -379.1428000000001
$ws.Cells.Item(86, 14).Value = -3382

# row 89 (hunk 17)
$ws.Cells.Item(89, 8).Value = 1392.3
$ws.Cells.Item(89, 9).Value = 1502.1428
$ws.Cells.Item(89, 10).Value = 1136
$ws.Cells.Item(89, 11).Value = 7510.714
$ws.Cells.Item(89, 12).Value = 5680
$ws.Cells.Item(89, 13).Value = -1894.714
$ws.Cells.Item(89, 14).Value = -16912

# row 100 (hunk 18)
$ws.Cells.Item(100, 8).Value = 18071.5
$ws.Cells.Item(100, 10).Value = 18071.5
$ws.Cells.Item(100, 12).Value = 18071.5
$ws.Cells.Item(100, 14).Value = -20235.5

# row 105 (hunk 19)
$ws.Cells.Item(105, 8).Value = 51087
$ws.Cells.Item(105, 9).Value = 1450
$ws.Cells.Item(105, 11).Value = 1450
$ws.Cells.Item(105, 13).Value = 297

# row 134 (hunk 20)
$ws.Cells.Item(134, 8).Value = 6205.5557
$ws.Cells.Item(134, 9).Value = 1462.5
$ws.Cells.Item(134, 11).Value = 4387.5
$ws.Cells.Item(134, 13).Value = -1852.5


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 8 (hunk 21)
$ws.Cells.Item(8, 8).Value = 600
$ws.Cells.Item(8, 10).Value = 600
$ws.Cells.Item(8, 12).Value = 600
$ws.Cells.Item(8, 14).Value = -880

# row 15 (hunk 22)
$ws.Cells.Item(15, 8).Value = 12736
$ws.Cells.Item(15, 10).Value = 12736
$ws.Cells.Item(15, 12).Value = 12736
$ws.Cells.Item(15, 14).Value = -13076

# row 62 (hunk 23)
$ws.Cells.Item(62, 8).Value = 4677.6
$ws.Cells.Item(62, 9).Value = 4677.6
$ws.Cells.Item(62, 11).Value = 4677.6
$ws.Cells.Item(62, 13).Value = -4053.6

# row 65 (hunk 24)
$ws.Cells.Item(65, 8).Value = 4677.6
$ws.Cells.Item(65, 9).Value = 4677.6
$ws.Cells.Item(65, 11).Value = 23388
$ws.Cells.Item(65, 13).Value = -20268


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 34 (hunk 25)
$ws.Cells.Item(34, 8).Value = 858.3333
$ws.Cells.Item(34, 10).Value = 1416.6666
$ws.Cells.Item(34, 12).Value = 4249.9998
$ws.Cells.Item(34, 14).Value = -4417.9998


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 3 (hunk 26)
$ws.Cells.Item(3, 8).Value = 8500000
$ws.Cells.Item(3, 9).Value = 12000000
$ws.Cells.Item(3, 11).Value = 12000000
$ws.Cells.Item(3, 13).Value = -11999884

# row 70 (hunk 27)
$ws.Cells.Item(70, 8).Value = 6000
$ws.Cells.Item(70, 10).Value = 6500
$ws.Cells.Item(70, 12).Value = 6500
$ws.Cells.Item(70, 14).Value = -7040

# row 73 (hunk 28)
$ws.Cells.Item(73, 8).Value = 6000
$ws.Cells.Item(73, 10).Value = 6500
$ws.Cells.Item(73, 12).Value = 6500
$ws.Cells.Item(73, 14).Value = -8372

# row 80 (hunk 29)
$ws.Cells.Item(80, 8).Value = 12927.875
$ws.Cells.Item(80, 9).Value = 11581.25
$ws.Cells.Item(80, 11).Value = 11581.25
$ws.Cells.Item(80, 13).Value = -10583.25

# row 83 (hunk 30)
$ws.Cells.Item(83, 8).Value = 12927.875
$ws.Cells.Item(83, 9).Value = 11581.25
$ws.Cells.Item(83, 11).Value = 57906.25
$ws.Cells.Item(83, 13).Value = -52914.25

# row 92 (hunk 31)
$ws.Cells.Item(92, 8).Value = 14699.667
$ws.Cells.Item(92, 10).Value = 14699.667
$ws.Cells.Item(92, 12).Value = 14699.667
$ws.Cells.Item(92, 14).Value = -18443.667

# row 102 (hunk 32)
$ws.Cells.Item(102, 8).Value = 344.25
$ws.Cells.Item(102, 9).Value = 350.57144
$ws.Cells.Item(102, 11).Value = 350.57144
$ws.Cells.Item(102, 13).Value = 1271.42856

# row 132 (hunk 33)
$ws.Cells.Item(132, 8).Value = 1738.5555
$ws.Cells.Item(132, 9).Value = 1492.4286
$ws.Cells.Item(132, 10).Value = 2600
$ws.Cells.Item(132, 11).Value = 4477.2858
$ws.Cells.Item(132, 12).Value = 7800
$ws.Cells.Item(132, 13).Value = -1947.2858
$ws.Cells.Item(132, 14).Value = -12860


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 29 (hunk 34)
$ws.Cells.Item(29, 8).Value = 28800
$ws.Cells.Item(29, 10).Value = 28800
$ws.Cells.Item(29, 12).Value = 28800
$ws.Cells.Item(29, 14).Value = -29390

# row 31 (hunk 35)
$ws.Cells.Item(31, 8).Value = 2000
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).ClearContents()

# row 34 (hunk 36)
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).ClearContents()

# row 61 (hunk 37)
$ws.Cells.Item(61, 8).Value = 1495
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 1495
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 1495
$ws.Cells.Item(61, 14).Value = -1899
$ws.Cells.Item(61, 13).ClearContents()

# row 100 (hunk 38)
$ws.Cells.Item(100, 8).Value = 1450
$ws.Cells.Item(100, 10).Value = 1450
$ws.Cells.Item(100, 12).Value = 1450
$ws.Cells.Item(100, 14).Value = -2532

# row 113 (hunk 39)
$ws.Cells.Item(113, 8).Value = 1495
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 1495
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 1495
$ws.Cells.Item(113, 14).Value = -5835
$ws.Cells.Item(113, 13).ClearContents()


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 92 (hunk 40)
$ws.Cells.Item(92, 8).Value = 35000
$ws.Cells.Item(92, 10).Value = 35000
$ws.Cells.Item(92, 12).Value = 35000
$ws.Cells.Item(92, 14).Value = -39992

# row 100 (hunk 41)
$ws.Cells.Item(100, 8).Value = 999
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()

# row 132 (hunk 42)
$ws.Cells.Item(132, 8).Value = 926
$ws.Cells.Item(132, 9).Value = 907.5
$ws.Cells.Item(132, 11).Value = 2722.5
$ws.Cells.Item(132, 13).Value = -192.5

